# Updated symbol list on Sat Dec 17 02:48:26 UTC 2022 with GitHub Actions
#
# Applies the latest crypto price snapshot to the "Price" (D) column and
# fixes a couple of stray "Worstin24h" suffixes that had leaked into the
# "Volume(1h)" (E) column text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values must be written as Text so Excel doesn't coerce
# them into floating point numbers (which would lose the original
# formatting/precision of the source string).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2")  "226.77"
Set-TextValue $ws.Range("D3")  "22.45"
Set-TextValue $ws.Range("D4")  "5.283"
Set-TextValue $ws.Range("D5")  "0.05531"
Set-TextValue $ws.Range("D6")  "3.386"
Set-TextValue $ws.Range("D7")  "6.473"
Set-TextValue $ws.Range("D8")  "0.7814"

Set-TextValue $ws.Range("D9")  "1.033"
$ws.Range("E9").Value = "8FTXTokenFTT"

Set-TextValue $ws.Range("D10") "0.1381"
Set-TextValue $ws.Range("D11") "0.07500"
Set-TextValue $ws.Range("D12") "0.03137"
Set-TextValue $ws.Range("D13") "0.02948"
Set-TextValue $ws.Range("D14") "0.09250"
Set-TextValue $ws.Range("D15") "0.001664"
Set-TextValue $ws.Range("D17") "0.04788"

Set-TextValue $ws.Range("D18") "0.0005864"
$ws.Range("E18").Value = "17OneONEWorstin24h"

Set-TextValue $ws.Range("D19") "0.006214"
Set-TextValue $ws.Range("D20") "0.005221"
Set-TextValue $ws.Range("D23") "3.833"
Set-TextValue $ws.Range("D26") "0.1286"
Set-TextValue $ws.Range("D27") "0.0005023"
Set-TextValue $ws.Range("D40") "0.03892"
Set-TextValue $ws.Range("D41") "0.007156"
Set-TextValue $ws.Range("D42") "0.1031"
Set-TextValue $ws.Range("D43") "0.003276"
Set-TextValue $ws.Range("D44") "0.009245"
Set-TextValue $ws.Range("D45") "0.00005426"
Set-TextValue $ws.Range("D46") "0.00000000751"
Set-TextValue $ws.Range("D47") "0.6757"
Set-TextValue $ws.Range("D48") "0.08882"
Set-TextValue $ws.Range("D49") "0.00002102"
Set-TextValue $ws.Range("D50") "0.01011"
